$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.800.05"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.854.10"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'313.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.4297"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").Value = "'0.3650"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").Value = "'44.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("D10").Value = "'0.07326"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").Value = "'0.8792"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").Value = "'20.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "1.870.64"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "'5.346"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").Value = "'6.533"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "'0.06908"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "'80.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.55%  "
$ws.Range("D19").Value = "'0.000009008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'15.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "27.787.53"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").Value = "'4.946"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("E24").Value = "  -2.65%  "
$ws.Range("D25").Value = "2.145.40"
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("D26").Value = "'1.991"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.91%  "
$ws.Range("D27").Value = "'155.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").Value = "'18.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.54%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'5.342"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'121.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.11%  "
$ws.Range("D31").Value = "'1.848"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("D32").Value = "'0.08906"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'0.7625"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").Value = "'4.554"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'1.108"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.61%  "
$ws.Range("D37").Value = "'0.05413"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "'1.100"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").Value = "'0.01945"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").Value = "'2.830"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.01%  "
$ws.Range("D41").Value = "'0.5080"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").Value = "'0.1657"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Value = "'6.684"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("D44").Value = "'8.370"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").Value = "'0.06544"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("D46").Value = "'10.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").Value = "'0.4678"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("D48").Value = "'104.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("D49").Value = "'0.9991"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "'1.622"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").Value = "'64.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.25%  "
